$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update data row 2
$ws.Range("B2").Value = 4.3508091582938189
$ws.Range("C2").Value = 6.317758840049482
$ws.Range("D2").Value = 3.817409235023514
$ws.Range("E2").Value = 3.9410391808256309

# Update data row 3
$ws.Range("B3").Value = 7.900352088866569
$ws.Range("C3").Value = 13.637047579725522
$ws.Range("D3").Value = 10.676319606241041
$ws.Range("E3").Value = 3.052185905350814

# Update the selection shown on the sheet to reflect the smaller range of interest
$excel.Goto($ws.Range("B1:E3"))
